$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 149.38
$ws.Range("F2").Value = 24237.06
$ws.Range("G2").Value = 155.68
$ws.Range("E4").Value = 16.28
$ws.Range("F4").Value = 411.96
$ws.Range("G4").Value = 20.3
$ws.Range("E5").Value = 150.73
$ws.Range("F5").Value = 24673.77
$ws.Range("G5").Value = 157.08
$ws.Range("E7").Value = 16.21
$ws.Range("F7").Value = 409.55
$ws.Range("G7").Value = 20.24
$ws.Range("E8").Value = 150.68
$ws.Range("F8").Value = 24659.81
$ws.Range("G8").Value = 157.03
$ws.Range("E10").Value = 16.24
$ws.Range("F10").Value = 409
$ws.Range("G10").Value = 20.22
$ws.Range("E11").Value = 30.21
$ws.Range("F11").Value = 1403.75
$ws.Range("G11").Value = 37.47
$ws.Range("E13").Value = 26.3
$ws.Range("F13").Value = 1021
$ws.Range("G13").Value = 31.95
$ws.Range("E14").Value = 30.28
$ws.Range("F14").Value = 1415.19
$ws.Range("G14").Value = 37.62
$ws.Range("E16").Value = 26.19
$ws.Range("F16").Value = 1060.89
$ws.Range("G16").Value = 32.57
$ws.Range("E17").Value = 30.34
$ws.Range("F17").Value = 1426.6
$ws.Range("G17").Value = 37.77
$ws.Range("E19").Value = 25.69
$ws.Range("F19").Value = 1012.45
$ws.Range("G19").Value = 31.82
$ws.Range("E20").Value = 4.85
$ws.Range("F20").Value = 34.2
$ws.Range("G20").Value = 5.85
$ws.Range("F22").Value = 26.48
$ws.Range("G22").Value = 5.15
$ws.Range("F23").Value = 25.75
$ws.Range("G23").Value = 5.07
$ws.Range("F25").Value = 28.64
$ws.Range("G25").Value = 5.35
$ws.Range("E26").Value = 4.83
$ws.Range("F26").Value = 25.8
$ws.Range("G26").Value = 5.08
$ws.Range("F28").Value = 31.61
$ws.Range("G28").Value = 5.62
